$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "AutomationManickam@gmail.com"
$ws.Range("C3").Value = "AutomationMonika@gmail.com"
$ws.Range("C4").Value = "AutomationHaritha@gmail.com"
